# Updates 후성.xlsx (IFRS financial figures, rows 2-9 / columns D:AJ)
# to the restated values from the commit "error solve ifrs list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1870
$ws.Range("E2").Value = -92
$ws.Range("F2").Value = -92
$ws.Range("G2").Value = -587
$ws.Range("H2").Value = -687
$ws.Range("I2").Value = -682
$ws.Range("J2").Value = -5
$ws.Range("K2").Value = 2320
$ws.Range("L2").Value = 1706
$ws.Range("M2").Value = 614
$ws.Range("N2").Value = 495
$ws.Range("O2").Value = 119
$ws.Range("P2").Value = 424
$ws.Range("Q2").Value = 79
$ws.Range("R2").Value = -137
$ws.Range("S2").Value = 70
$ws.Range("T2").Value = 198
$ws.Range("U2").Value = -120
$ws.Range("V2").Value = 1310
$ws.Range("W2").Value = -4.93
$ws.Range("X2").Value = -36.73
$ws.Range("Y2").Value = -80.58
$ws.Range("Z2").Value = -26.31
$ws.Range("AA2").Value = 277.99
$ws.Range("AB2").Value = 10.27
$ws.Range("AC2").Value = -784
$ws.Range("AD2").Value = -3.98
$ws.Range("AE2").Value = 571
$ws.Range("AF2").Value = 5.47
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 86918416

# Row 3
$ws.Range("D3").Value = 1612
$ws.Range("E3").Value = 152
$ws.Range("F3").Value = 156
$ws.Range("G3").Value = 131
$ws.Range("H3").Value = 137
$ws.Range("I3").Value = 139
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 2614
$ws.Range("L3").Value = 1634
$ws.Range("M3").Value = 980
$ws.Range("N3").Value = 861
$ws.Range("O3").Value = 119
$ws.Range("P3").Value = 462
$ws.Range("Q3").Value = 346
$ws.Range("R3").Value = -409
$ws.Range("S3").Value = 117
$ws.Range("T3").Value = 363
$ws.Range("U3").Value = -16
$ws.Range("V3").Value = 1215
$ws.Range("W3").Value = 9.44
$ws.Range("X3").Value = 8.51
$ws.Range("Y3").Value = 20.56
$ws.Range("Z3").Value = 5.56
$ws.Range("AA3").Value = 166.74
$ws.Range("AB3").Value = 75.84999999999999
$ws.Range("AC3").Value = 157
$ws.Range("AD3").Value = 34.06
$ws.Range("AE3").Value = 935
$ws.Range("AF3").Value = 5.73
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 92399753

# Row 4
$ws.Range("D4").Value = 1921
$ws.Range("E4").Value = 364
$ws.Range("F4").Value = 364
$ws.Range("G4").Value = 594
$ws.Range("H4").Value = 614
$ws.Range("I4").Value = 609
$ws.Range("J4").Value = 6
$ws.Range("K4").Value = 3109
$ws.Range("L4").Value = 1517
$ws.Range("M4").Value = 1592
$ws.Range("N4").Value = 1472
$ws.Range("O4").Value = 120
$ws.Range("P4").Value = 462
$ws.Range("Q4").Value = 349
$ws.Range("R4").Value = -133
$ws.Range("S4").Value = -68
$ws.Range("T4").Value = 98
$ws.Range("U4").Value = 251
$ws.Range("V4").Value = 1024
$ws.Range("W4").Value = 18.97
$ws.Range("X4").Value = 31.99
$ws.Range("Y4").Value = 52.17
$ws.Range("Z4").Value = 21.47
$ws.Range("AA4").Value = 95.28
$ws.Range("AB4").Value = 207.43
$ws.Range("AC4").Value = 659
$ws.Range("AD4").Value = 11.14
$ws.Range("AE4").Value = 1595
$ws.Range("AF4").Value = 4.6
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 92399753

# Row 5
$ws.Range("D5").Value = 2488
$ws.Range("E5").Value = 355
$ws.Range("F5").Value = 355
$ws.Range("G5").Value = 333
$ws.Range("H5").Value = 282
$ws.Range("I5").Value = 295
$ws.Range("J5").Value = -13
$ws.Range("K5").Value = 3161
$ws.Range("L5").Value = 1072
$ws.Range("M5").Value = 2089
$ws.Range("N5").Value = 1765
$ws.Range("O5").Value = 325
$ws.Range("P5").Value = 462
$ws.Range("Q5").Value = 519
$ws.Range("R5").Value = -316
$ws.Range("S5").Value = -42
$ws.Range("T5").Value = 296
$ws.Range("U5").Value = 223
$ws.Range("V5").Value = 744
$ws.Range("W5").Value = 14.26
$ws.Range("X5").Value = 11.34
$ws.Range("Y5").Value = 18.22
$ws.Range("Z5").Value = 9
$ws.Range("AA5").Value = 51.31
$ws.Range("AB5").Value = 271.38
$ws.Range("AC5").Value = 319
$ws.Range("AD5").Value = 31.8
$ws.Range("AE5").Value = 1912
$ws.Range("AF5").Value = 5.31
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 92399753

# Row 6
$ws.Range("D6").Value = 2749
$ws.Range("E6").Value = 396
$ws.Range("F6").Value = 396
$ws.Range("G6").Value = 359
$ws.Range("H6").Value = 281
$ws.Range("I6").Value = 308
$ws.Range("K6").Value = 4415
$ws.Range("L6").Value = 1960
$ws.Range("M6").Value = 2455
$ws.Range("N6").Value = 2080
$ws.Range("P6").Value = 463
$ws.Range("Q6").Value = 307
$ws.Range("R6").Value = -1182
$ws.Range("S6").Value = 850
$ws.Range("T6").Value = 1223
$ws.Range("U6").Value = -916
$ws.Range("V6").Value = 1491
$ws.Range("W6").Value = 14.4
$ws.Range("X6").Value = 10.2
$ws.Range("Y6").Value = 16.04
$ws.Range("Z6").Value = 7.41
$ws.Range("AA6").Value = 79.83
$ws.Range("AB6").Value = 339.43
$ws.Range("AC6").Value = 333
$ws.Range("AD6").Value = 22.27
$ws.Range("AE6").Value = 2246
$ws.Range("AF6").Value = 3.3
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 92606819

# Row 7
$ws.Range("D7").Value = 2579
$ws.Range("E7").Value = 248
$ws.Range("G7").Value = 216
$ws.Range("H7").Value = 177
$ws.Range("I7").Value = 201
$ws.Range("K7").Value = 5631
$ws.Range("L7").Value = 2974
$ws.Range("M7").Value = 2656
$ws.Range("N7").Value = 2303
$ws.Range("P7").Value = 463
$ws.Range("Q7").Value = 538
$ws.Range("R7").Value = -1051
$ws.Range("S7").Value = 863
$ws.Range("T7").Value = 1070
$ws.Range("U7").Value = -450
$ws.Range("W7").Value = 9.6
$ws.Range("X7").Value = 6.86
$ws.Range("Y7").Value = 9.19
$ws.Range("Z7").Value = 3.52
$ws.Range("AA7").Value = 111.97
$ws.Range("AC7").Value = 217
$ws.Range("AD7").Value = 38.73
$ws.Range("AE7").Value = 2487
$ws.Range("AF7").Value = 3.39
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 2988
$ws.Range("E8").Value = 442
$ws.Range("G8").Value = 438
$ws.Range("H8").Value = 360
$ws.Range("I8").Value = 388
$ws.Range("K8").Value = 5978
$ws.Range("L8").Value = 2918
$ws.Range("M8").Value = 3060
$ws.Range("N8").Value = 2720
$ws.Range("P8").Value = 463
$ws.Range("Q8").Value = 648
$ws.Range("R8").Value = -282
$ws.Range("S8").Value = -115
$ws.Range("T8").Value = 278
$ws.Range("U8").Value = 366
$ws.Range("W8").Value = 14.8
$ws.Range("X8").Value = 12.06
$ws.Range("Y8").Value = 15.46
$ws.Range("Z8").Value = 6.21
$ws.Range("AA8").Value = 95.34
$ws.Range("AC8").Value = 419
$ws.Range("AD8").Value = 20.08
$ws.Range("AE8").Value = 2937
$ws.Range("AF8").Value = 2.87
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 3584
$ws.Range("E9").Value = 590
$ws.Range("G9").Value = 660
$ws.Range("H9").Value = 561
$ws.Range("I9").Value = 488
$ws.Range("K9").Value = 6768
$ws.Range("L9").Value = 3169
$ws.Range("M9").Value = 3599
$ws.Range("N9").ClearContents()
$ws.Range("P9").Value = 463
$ws.Range("Q9").Value = 806
$ws.Range("R9").Value = -605
$ws.Range("S9").Value = -100
$ws.Range("T9").Value = 600
$ws.Range("U9").Value = 201
$ws.Range("W9").Value = 16.48
$ws.Range("X9").Value = 15.65
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").Value = 8.800000000000001
$ws.Range("AA9").Value = 88.05
$ws.Range("AC9").Value = 527
$ws.Range("AD9").Value = 15.98
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
